# Update the weekly Fruta/Hortaliza price records (rows 2-7) on the active sheet.
# The underlying source data was refreshed; several rows now carry different
# dates, volumes, prices, units, origins and quality grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44334
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = "$/caja 12 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 11500
$ws.Range("T2").Value = 1

# Row 3
$ws.Range("D3").Value = 44330
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15500
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 861

# Row 4
$ws.Range("D4").Value = 44742
$ws.Range("L4").Value = "Segunda"
$ws.Range("P4").Value = 14500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 806

# Row 5
$ws.Range("D5").Value = 44708
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12571
$ws.Range("Q5").Value = "$/caja 12 kilos empedrada"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 1048
$ws.Range("T5").Value = 12

# Row 6
$ws.Range("D6").Value = 44714
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 806

# Row 7
$ws.Range("D7").Value = 44719
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("P7").Value = 14400
$ws.Range("R7").Value = "Región del Maule"
$ws.Range("S7").Value = 800
